$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.709.76"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "1.545.08"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.15%  "
$ws.Range("E9").Value = "  -1.18%  "
$ws.Range("E10").Value = "  -0.50%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "1.764.32"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").Value = "1.538.80"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.509"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "26.704.45"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "0.0₃0688"
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  -4.20%  "
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0458"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("D33").Value = "1.347.27"
$ws.Range("E33").Value = "  -2.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  -2.93%  "
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.934"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.800"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.35%  "
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").Value = "1.678.35"
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0508"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "0.0₇0973"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  +0.51%  "
